$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C52").Value = "linktest"
$ws.Range("C52").Font.Color = 16711680
$ws.Range("C52").Font.Underline = $true
$ws.Range("C52").Font.Name = "Calibri"
